$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 8:9 get their own (new) cell style distinct from the one
# used by rows 2-7: explicitly pin the font to the theme's minor font
# with solid black, which mints a fresh font + cellXfs entry.
$ws.Range("A8:L9").Font.Color = 0
$ws.Range("A8:L9").Font.ThemeFont = 1

$ws.Range("A8").Value = "lab_4821"
$ws.Range("B8").Value = "Reddy_4821_180419A3"
$ws.Range("C8").Value = "XXXXXXX"
$ws.Range("D8").Value = "Sample_3"
$ws.Range("F8").Value = "PE"
$ws.Range("G8").Value = "hg38"
$ws.Range("H8").Value = "STARR-seq"
$ws.Range("E8").Value = "CAL51.starrseq.Gefit.rep1"

$ws.Range("A9").Value = "lab_4821"
$ws.Range("B9").Value = "Reddy_4821_180419A3"
$ws.Range("C9").Value = "XXXXXXX"
$ws.Range("D9").Value = "Sample_4"
$ws.Range("F9").Value = "PE"
$ws.Range("G9").Value = "hg38"
$ws.Range("H9").Value = "STARR-seq"
$ws.Range("E9").Value = "CAL51.starrseq.Gefit_Inibit.rep1"

$ws.Range("F17").Select() | Out-Null
